# IST price update 2025-12-23 09:02
# A new price-scrape column is inserted right after the "SKU Name" column
# (i.e. before the current column B), pushing every existing timestamp
# column one position to the right (B->C, C->D, ... AH->AI) and growing
# the used range from A1:AH26 to A1:AI26.
#
# The brand-new column gets the latest scrape timestamp in row 1, and for
# every product row it is seeded with that product's most-recent known
# price (i.e. what used to be in column B, now shifted into column C) -
# except for row 13, where this scrape failed to capture a price, so the
# new cell is left blank, matching the source price tracker's behaviour
# for missing data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B; this shifts B:AH -> C:AI and
# extends the sheet dimension to A1:AI26 automatically.
$ws.Columns.Item(2).Insert()

# Stamp the new column's header with the new scrape timestamp.
$ws.Cells.Item(1, 2).Value2 = "2025-12-23 14:27"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 13) {
        # This product's price was not captured in this scrape run -
        # leave the new column blank for it.
        continue
    }
    # Carry the previous latest price (now in column C, post-insert)
    # forward into the brand-new column B.
    $ws.Cells.Item($r, 2).Value2 = $ws.Cells.Item($r, 3).Value2
}
